$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (11) down into the two new rows
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D13").PasteSpecial(-4122)

# Row 12: .about-us .content-left .title
$ws.Range("A12").Value = ".about-us .content-left .title"
$ws.Range("B12").Value = "font-size"
$ws.Range("C12").Value = 35.11
$ws.Range("D12").Formula = "=C12/B1"

# Row 13: .about-us .content-left .content-text
$ws.Range("A13").Value = ".about-us .content-left .content-text"
$ws.Range("B13").Value = "font-size"
$ws.Range("C13").Value = 18.06
$ws.Range("D13").Formula = "=C13/B1"

$ws.Range("E15").Select()
